# Update "想去人数" (want-to-go headcount) figures pulled from bilibili show
# pages for the 南宁-漫展信息 workbook (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 66    # 南宁·三月三漫次元国风动漫节: 59 -> 66
$wsExpo.Range("F4").Value = 2329  # 南宁·2024三月三国潮动漫节（良牙春典）: 2313 -> 2329
$wsExpo.Range("F6").Value = 511   # 南宁·布谷鸟动漫展4th: 508 -> 511

# Sheet "全部类型" (all types, combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 66     # 南宁·三月三漫次元国风动漫节: 59 -> 66
$wsAll.Range("F6").Value = 2330   # 南宁·2024三月三国潮动漫节（良牙春典）: 2313 -> 2330
$wsAll.Range("F8").Value = 511    # 南宁·布谷鸟动漫展4th: 508 -> 511
